$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misread as a number (e.g. "0.9997") are
# written via a Text-format round trip so they stay text, matching the
# original inlineStr cells, then restored to the default "Normal" style
# so no stray cell formatting is introduced.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '30.748.57'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '1.890.96'
$ws.Range('E3').Value = '  +0.36%  '
Set-TextValue 'D4' '0.9997'
$ws.Range('E4').Value = '  -0.10%  '
Set-TextValue 'D5' '249.61'
$ws.Range('E5').Value = '  +1.01%  '
Set-TextValue 'D6' '0.9996'
$ws.Range('E6').Value = '  -0.05%  '
Set-TextValue 'D7' '0.4764'
$ws.Range('E7').Value = '  -0.15%  '
Set-TextValue 'D8' '0.2937'
$ws.Range('E8').Value = '  +0.47%  '
Set-TextValue 'D9' '0.06542'
$ws.Range('E9').Value = '  +0.10%  '
Set-TextValue 'D10' '22.11'
$ws.Range('E10').Value = '  +0.27%  '
Set-TextValue 'D11' '0.07758'
$ws.Range('E11').Value = '  +0.40%  '
Set-TextValue 'D12' '97.27'
$ws.Range('E12').Value = '  -0.45%  '
Set-TextValue 'D13' '0.7404'
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('D14').Value = '1.888.49'
$ws.Range('E14').Value = '  +0.25%  '
Set-TextValue 'D15' '5.254'
$ws.Range('E15').Value = '  +1.76%  '
Set-TextValue 'D16' '283.84'
$ws.Range('E16').Value = '  +3.10%  '
$ws.Range('D17').Value = '30.795.81'
$ws.Range('E17').Value = '  +0.70%  '
Set-TextValue 'D18' '13.21'
$ws.Range('E18').Value = '  -2.44%  '
Set-TextValue 'D19' '0.000007581'
$ws.Range('E19').Value = '  -0.10%  '
Set-TextValue 'D20' '0.9997'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '2.138.06'
$ws.Range('E21').Value = '  +0.62%  '
Set-TextValue 'D22' '5.342'
$ws.Range('E22').Value = '  +1.25%  '
Set-TextValue 'D23' '0.9990'
$ws.Range('E23').Value = '  -0.18%  '
Set-TextValue 'D24' '6.257'
$ws.Range('E24').Value = '  +0.64%  '
Set-TextValue 'D25' '9.259'
$ws.Range('E25').Value = '  -0.93%  '
Set-TextValue 'D26' '164.39'
$ws.Range('E26').Value = '  +0.43%  '
Set-TextValue 'D27' '18.95'
$ws.Range('E27').Value = '  -0.05%  '
Set-TextValue 'D28' '1.932'
$ws.Range('E28').Value = '  -0.92%  '
Set-TextValue 'D29' '1.345'
$ws.Range('E29').Value = '  -2.11%  '
Set-TextValue 'D30' '0.09754'
$ws.Range('E30').Value = '  -2.38%  '
Set-TextValue 'D31' '1.503'
$ws.Range('E31').Value = '  -0.86%  '
Set-TextValue 'D32' '4.315'
$ws.Range('E32').Value = '  -0.38%  '
Set-TextValue 'D33' '4.210'
$ws.Range('E33').Value = '  +2.35%  '
Set-TextValue 'D34' '0.04895'
$ws.Range('E34').Value = '  +1.73%  '
Set-TextValue 'D35' '1.129'
$ws.Range('E35').Value = '  -0.16%  '
Set-TextValue 'D36' '0.7014'
$ws.Range('E36').Value = '  -0.48%  '
Set-TextValue 'D37' '2.721'
$ws.Range('E37').Value = '  +0.08%  '
Set-TextValue 'D38' '0.01919'
$ws.Range('E38').Value = '  +2.46%  '
Set-TextValue 'D39' '2.807'
$ws.Range('E39').Value = '  +2.06%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D40' '6.354'
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D41' '76.27'
$ws.Range('E41').Value = '  +6.77%  '
Set-TextValue 'D42' '2.033'
$ws.Range('E42').Value = '  +3.00%  '
Set-TextValue 'D43' '0.4273'
$ws.Range('E43').Value = '  +1.09%  '
Set-TextValue 'D44' '0.8423'
$ws.Range('E44').Value = '  +0.18%  '
Set-TextValue 'D45' '0.9997'
$ws.Range('E45').Value = '  -0.01%  '
Set-TextValue 'D46' '101.97'
$ws.Range('E46').Value = '  -0.93%  '
Set-TextValue 'D47' '9.449'
$ws.Range('E47').Value = '  +1.57%  '
Set-TextValue 'D48' '7.100'
$ws.Range('E48').Value = '  -0.20%  '
Set-TextValue 'D49' '35.81'
$ws.Range('E49').Value = '  +0.36%  '
Set-TextValue 'D50' '924.79'
$ws.Range('E50').Value = '  +0.77%  '
Set-TextValue 'D51' '0.05778'
$ws.Range('E51').Value = '  +2.33%  '
